# Insert a new data row at row 205 (pushes existing rows 205..290 down to
# 206..291) and populate it with the new Rabanito/Vega Central record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 205..290 down by inserting a new blank row at 205.
$ws.Rows.Item(205).Insert()

# Populate the newly inserted row 205.
$ws.Cells.Item(205, 1).Value = 9
$ws.Cells.Item(205, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(205, 3).Value = "Metropolitana"
$ws.Cells.Item(205, 4).Value = 44755
$ws.Cells.Item(205, 5).Value = 13
$ws.Cells.Item(205, 6).Value = 300000001
$ws.Cells.Item(205, 7).Value = "Rabanito"
$ws.Cells.Item(205, 8).Value = "Sin especificar"
$ws.Cells.Item(205, 9).Value = "Primera"
$ws.Cells.Item(205, 10).Value = 7000
$ws.Cells.Item(205, 11).Value = 2500
$ws.Cells.Item(205, 12).Value = 3000
$ws.Cells.Item(205, 13).Value = 2750
$ws.Cells.Item(205, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(205, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(205, 16).Value = 28
$ws.Cells.Item(205, 17).Value = 100
$ws.Cells.Item(205, 18).Value = "Hortaliza"
